$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value2 = '[{"role": "assistant", "content": "Chào cậu! Bữa nay thầy giáo robot trên sao HỎA dạy tớ cách luyện nói theo nhịp như… hát rap đó! Mỗi cụm từ giống như một đoạn nhạc nhỏ, có nhịp điệu riêng."}, {"role": "assistant", "content": "Cứ theo nhịp là nói được tiếng Anh hay hơn liền! Giờ tụi mình thử cùng luyện nha. Cậu đã sẵn sàng chưa?"}, {"role": "user", "content": "[NEXT]"}]'
$ws.Range("C2").Value2 = 'Chào cậu! Bữa nay thầy giáo robot trên sao HỎA dạy tớ cách luyện nói theo nhịp như… hát rap đó! Mỗi cụm từ giống như một đoạn nhạc nhỏ, có nhịp điệu riêng.'

$ws.Range("A3").Value2 = '[{"role": "assistant", "content": "Chào cậu! Bữa nay thầy giáo robot trên sao HỎA dạy tớ cách luyện nói theo nhịp như… hát rap đó! Mỗi cụm từ giống như một đoạn nhạc nhỏ, có nhịp điệu riêng."}, {"role": "assistant", "content": "Cứ theo nhịp là nói được tiếng Anh hay hơn liền! Giờ tụi mình thử cùng luyện nha. Cậu đã sẵn sàng chưa?"}, {"role": "user", "content": "[NEXT]"}, {"role": "assistant", "content": "Yên tâm nha! Cậu không học một mình đâu, Pika luôn bên cạnh nè! chúng mình sẽ học nói tiếng Anh theo nhịp giống như hát nha! Bài hát này có nhiều cụm từ rất quen mà cũng rất vui. Bài đầu tiên nha. Đầu tiên, cậu hãy nghe thử cả bài nhé, đừng hát theo vội. Mình chỉ cần cảm nhịp thôi."}, {"role": "assistant", "content": "Giờ mình cùng lặp lại đoạn đầu tiên nha, câu chào cực kỳ thân thiện luôn nè."}, {"role": "assistant", "content": "Cậu nói lại cùng tớ nha."}, {"role": "user", "content": "Hello, hello, hello, how are you? Hello, hello, hello, how are you?"}]'
$ws.Range("C3").Value2 = 'Cứ theo nhịp là nói được tiếng Anh hay hơn liền! Giờ tụi mình thử cùng luyện nha. Cậu đã sẵn sàng chưa?'

$ws.Range("A4").Value2 = '[{"role": "assistant", "content": "Chào cậu! Bữa nay thầy giáo robot trên sao HỎA dạy tớ cách luyện nói theo nhịp như… hát rap đó! Mỗi cụm từ giống như một đoạn nhạc nhỏ, có nhịp điệu riêng."}, {"role": "assistant", "content": "Cứ theo nhịp là nói được tiếng Anh hay hơn liền! Giờ tụi mình thử cùng luyện nha. Cậu đã sẵn sàng chưa?"}, {"role": "user", "content": "[NEXT]"}, {"role": "assistant", "content": "Yên tâm nha! Cậu không học một mình đâu, Pika luôn bên cạnh nè! chúng mình sẽ học nói tiếng Anh theo nhịp giống như hát nha! Bài hát này có nhiều cụm từ rất quen mà cũng rất vui. Bài đầu tiên nha. Đầu tiên, cậu hãy nghe thử cả bài nhé, đừng hát theo vội. Mình chỉ cần cảm nhịp thôi."}, {"role": "assistant", "content": "Giờ mình cùng lặp lại đoạn đầu tiên nha, câu chào cực kỳ thân thiện luôn nè."}, {"role": "assistant", "content": "Cậu nói lại cùng tớ nha."}, {"role": "user", "content": "Hello, hello, hello, how are you? Hello, hello, hello, how are you?"}, {"role": "assistant", "content": "Nice! You did it! Tiếp theo là khi mình cảm thấy thật tuyệt, cùng nói theo nhịp nha!"}, {"role": "assistant", "content": "Nào, nói lại với tớ!"}, {"role": "user", "content": "I''m good. I''m rich. I''m wonderful."}]'
$ws.Range("C4").Value2 = 'Yên tâm nha! Cậu không học một mình đâu, Pika luôn bên cạnh nè! chúng mình sẽ học nói tiếng Anh theo nhịp giống như hát nha! Bài hát này có nhiều cụm từ rất quen mà cũng rất vui. Bài đầu tiên nha. Đầu tiên, cậu hãy nghe thử cả bài nhé, đừng hát theo vội. Mình chỉ cần cảm nhịp thôi. '

$ws.Range("A5").Value2 = '[{"role": "assistant", "content": "Chào cậu! Bữa nay thầy giáo robot trên sao HỎA dạy tớ cách luyện nói theo nhịp như… hát rap đó! Mỗi cụm từ giống như một đoạn nhạc nhỏ, có nhịp điệu riêng."}, {"role": "assistant", "content": "Cứ theo nhịp là nói được tiếng Anh hay hơn liền! Giờ tụi mình thử cùng luyện nha. Cậu đã sẵn sàng chưa?"}, {"role": "user", "content": "[NEXT]"}, {"role": "assistant", "content": "Yên tâm nha! Cậu không học một mình đâu, Pika luôn bên cạnh nè! chúng mình sẽ học nói tiếng Anh theo nhịp giống như hát nha! Bài hát này có nhiều cụm từ rất quen mà cũng rất vui. Bài đầu tiên nha. Đầu tiên, cậu hãy nghe thử cả bài nhé, đừng hát theo vội. Mình chỉ cần cảm nhịp thôi."}, {"role": "assistant", "content": "Giờ mình cùng lặp lại đoạn đầu tiên nha, câu chào cực kỳ thân thiện luôn nè."}, {"role": "assistant", "content": "Cậu nói lại cùng tớ nha."}, {"role": "user", "content": "Hello, hello, hello, how are you? Hello, hello, hello, how are you?"}, {"role": "assistant", "content": "Nice! You did it! Tiếp theo là khi mình cảm thấy thật tuyệt, cùng nói theo nhịp nha!"}, {"role": "assistant", "content": "Nào, nói lại với tớ!"}, {"role": "user", "content": "I''m good. I''m rich. I''m wonderful."}, {"role": "assistant", "content": "Awesome! Cậu nói theo nhịp rất tốt! Câu chào quen thuộc quay lại nè! Mình cùng nói lại một lần nữa nhé!"}, {"role": "assistant", "content": "Cậu nhớ nhịp rồi đúng không? Cùng nói nào!"}, {"role": "user", "content": "Hello! Hello!"}]'
$ws.Range("C5").Value2 = ""

$ws.Range("A6").Value2 = '[{"role": "assistant", "content": "Chào cậu! Bữa nay thầy giáo robot trên sao HỎA dạy tớ cách luyện nói theo nhịp như… hát rap đó! Mỗi cụm từ giống như một đoạn nhạc nhỏ, có nhịp điệu riêng."}, {"role": "assistant", "content": "Cứ theo nhịp là nói được tiếng Anh hay hơn liền! Giờ tụi mình thử cùng luyện nha. Cậu đã sẵn sàng chưa?"}, {"role": "user", "content": "[NEXT]"}, {"role": "assistant", "content": "Yên tâm nha! Cậu không học một mình đâu, Pika luôn bên cạnh nè! chúng mình sẽ học nói tiếng Anh theo nhịp giống như hát nha! Bài hát này có nhiều cụm từ rất quen mà cũng rất vui. Bài đầu tiên nha. Đầu tiên, cậu hãy nghe thử cả bài nhé, đừng hát theo vội. Mình chỉ cần cảm nhịp thôi."}, {"role": "assistant", "content": "Giờ mình cùng lặp lại đoạn đầu tiên nha, câu chào cực kỳ thân thiện luôn nè."}, {"role": "assistant", "content": "Cậu nói lại cùng tớ nha."}, {"role": "user", "content": "Hello, hello, hello, how are you? Hello, hello, hello, how are you?"}, {"role": "assistant", "content": "Nice! You did it! Tiếp theo là khi mình cảm thấy thật tuyệt, cùng nói theo nhịp nha!"}, {"role": "assistant", "content": "Nào, nói lại với tớ!"}, {"role": "user", "content": "I''m good. I''m rich. I''m wonderful."}, {"role": "assistant", "content": "Awesome! Cậu nói theo nhịp rất tốt! Câu chào quen thuộc quay lại nè! Mình cùng nói lại một lần nữa nhé!"}, {"role": "assistant", "content": "Cậu nhớ nhịp rồi đúng không? Cùng nói nào!"}, {"role": "user", "content": "Hello! Hello!"}, {"role": "assistant", "content": "Good job! Cậu nhớ rất giỏi! Đôi khi mình cũng cảm thấy hơi mệt hoặc đói. Vậy thì mình sẽ nói thế này nè"}, {"role": "assistant", "content": "Giờ thì cậu nhại theo tớ nha!"}, {"role": "user", "content": "I''m tired I''m hungry I''m not so good I''m tired"}]'
$ws.Range("C6").Value2 = 'Cậu nói lại cùng tớ nha.'

$ws.Range("A7").Value2 = '[{"role": "assistant", "content": "Chào cậu! Bữa nay thầy giáo robot trên sao HỎA dạy tớ cách luyện nói theo nhịp như… hát rap đó! Mỗi cụm từ giống như một đoạn nhạc nhỏ, có nhịp điệu riêng."}, {"role": "assistant", "content": "Cứ theo nhịp là nói được tiếng Anh hay hơn liền! Giờ tụi mình thử cùng luyện nha. Cậu đã sẵn sàng chưa?"}, {"role": "user", "content": "[NEXT]"}, {"role": "assistant", "content": "Yên tâm nha! Cậu không học một mình đâu, Pika luôn bên cạnh nè! chúng mình sẽ học nói tiếng Anh theo nhịp giống như hát nha! Bài hát này có nhiều cụm từ rất quen mà cũng rất vui. Bài đầu tiên nha. Đầu tiên, cậu hãy nghe thử cả bài nhé, đừng hát theo vội. Mình chỉ cần cảm nhịp thôi."}, {"role": "assistant", "content": "Giờ mình cùng lặp lại đoạn đầu tiên nha, câu chào cực kỳ thân thiện luôn nè."}, {"role": "assistant", "content": "Cậu nói lại cùng tớ nha."}, {"role": "user", "content": "Hello, hello, hello, how are you? Hello, hello, hello, how are you?"}, {"role": "assistant", "content": "Nice! You did it! Tiếp theo là khi mình cảm thấy thật tuyệt, cùng nói theo nhịp nha!"}, {"role": "assistant", "content": "Nào, nói lại với tớ!"}, {"role": "user", "content": "I''m good. I''m rich. I''m wonderful."}, {"role": "assistant", "content": "Awesome! Cậu nói theo nhịp rất tốt! Câu chào quen thuộc quay lại nè! Mình cùng nói lại một lần nữa nhé!"}, {"role": "assistant", "content": "Cậu nhớ nhịp rồi đúng không? Cùng nói nào!"}, {"role": "user", "content": "Hello! Hello!"}, {"role": "assistant", "content": "Good job! Cậu nhớ rất giỏi! Đôi khi mình cũng cảm thấy hơi mệt hoặc đói. Vậy thì mình sẽ nói thế này nè"}, {"role": "assistant", "content": "Giờ thì cậu nhại theo tớ nha!"}, {"role": "user", "content": "I''m tired I''m hungry I''m not so good I''m tired"}, {"role": "assistant", "content": "Great! Bây giờ mình cùng nói câu chào cuối nhé! Và cuối cùng, bài hát kết thúc bằng câu chào quen thuộc một lần nữa. Nói lại lần cuối nha!"}, {"role": "assistant", "content": "Cùng nói lại với tớ nào!"}, {"role": "user", "content": "Hello, hello, hello, how are you?"}]'
$ws.Range("C7").Value2 = 'Nice! You did it! Tiếp theo là khi mình cảm thấy thật tuyệt, cùng nói theo nhịp nha! '

$ws.Range("A8").Value2 = '[{"role": "assistant", "content": "Chào cậu! Bữa nay thầy giáo robot trên sao HỎA dạy tớ cách luyện nói theo nhịp như… hát rap đó! Mỗi cụm từ giống như một đoạn nhạc nhỏ, có nhịp điệu riêng."}, {"role": "assistant", "content": "Cứ theo nhịp là nói được tiếng Anh hay hơn liền! Giờ tụi mình thử cùng luyện nha. Cậu đã sẵn sàng chưa?"}, {"role": "user", "content": "[NEXT]"}, {"role": "assistant", "content": "Yên tâm nha! Cậu không học một mình đâu, Pika luôn bên cạnh nè! chúng mình sẽ học nói tiếng Anh theo nhịp giống như hát nha! Bài hát này có nhiều cụm từ rất quen mà cũng rất vui. Bài đầu tiên nha. Đầu tiên, cậu hãy nghe thử cả bài nhé, đừng hát theo vội. Mình chỉ cần cảm nhịp thôi."}, {"role": "assistant", "content": "Giờ mình cùng lặp lại đoạn đầu tiên nha, câu chào cực kỳ thân thiện luôn nè."}, {"role": "assistant", "content": "Cậu nói lại cùng tớ nha."}, {"role": "user", "content": "Hello, hello, hello, how are you? Hello, hello, hello, how are you?"}, {"role": "assistant", "content": "Nice! You did it! Tiếp theo là khi mình cảm thấy thật tuyệt, cùng nói theo nhịp nha!"}, {"role": "assistant", "content": "Nào, nói lại với tớ!"}, {"role": "user", "content": "I''m good. I''m rich. I''m wonderful."}, {"role": "assistant", "content": "Awesome! Cậu nói theo nhịp rất tốt! Câu chào quen thuộc quay lại nè! Mình cùng nói lại một lần nữa nhé!"}, {"role": "assistant", "content": "Cậu nhớ nhịp rồi đúng không? Cùng nói nào!"}, {"role": "user", "content": "Hello! Hello!"}, {"role": "assistant", "content": "Good job! Cậu nhớ rất giỏi! Đôi khi mình cũng cảm thấy hơi mệt hoặc đói. Vậy thì mình sẽ nói thế này nè"}, {"role": "assistant", "content": "Giờ thì cậu nhại theo tớ nha!"}, {"role": "user", "content": "I''m tired I''m hungry I''m not so good I''m tired"}, {"role": "assistant", "content": "Great! Bây giờ mình cùng nói câu chào cuối nhé! Và cuối cùng, bài hát kết thúc bằng câu chào quen thuộc một lần nữa. Nói lại lần cuối nha!"}, {"role": "assistant", "content": "Cùng nói lại với tớ nào!"}, {"role": "user", "content": "Hello, hello, hello, how are you?"}, {"role": "assistant", "content": "Well done! Giờ thì mình cùng hát lại từ đầu đến cuối bài hát này nhé, lần này cậu có thể vừa hát vừa nhún nhảy cho vui luôn nha!"}, {"role": "assistant", "content": "Tiếp theo, hãy tưởng tượng mình sẽ bắt đầu một ngày mới thật vui bằng một bài hát tiếng Anh siêu rộn ràng nha. Cậu sẽ được luyện các câu chào buổi sáng theo nhịp. Vừa học, vừa hát, vừa nhảy cũng được luôn! Cậu sẵn sàng chưa?"}, {"role": "user", "content": "I''m ready."}]'
$ws.Range("C8").Value2 = 'Nào, nói lại với tớ!'

$ws.Rows("9:27").Delete()
